# Adds two new columns, I (I0) and J (IF), to the sheet, matching the
# header style used by the existing columns and filling in the per-row
# values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for the two new columns, styled like the existing header row
# (column H header "IP" carries the bold/bordered/centered header style).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Per-row data: row, I value, J value
$data = @(
    @(2, 7, 7),
    @(3, 1, 1),
    @(4, 6, 6),
    @(5, 6, 6),
    @(6, 1, 1),
    @(7, 5, 5),
    @(8, 7, 7),
    @(9, 5, 5),
    @(10, 7, 7),
    @(11, 1, 1),
    @(12, 4, 4),
    @(13, 1, 1),
    @(14, 9, 9),
    @(15, 7, 7),
    @(16, 1, 1),
    @(17, 9, 9),
    @(18, 7, 8),
    @(19, 9, 9),
    @(20, 8, 8),
    @(21, 7, 7),
    @(22, 6, 6),
    @(23, 5, 5),
    @(24, 7, 7),
    @(25, 7, 7),
    @(26, 5, 5),
    @(27, 6, 6),
    @(28, 8, 8),
    @(29, 6, 7),
    @(30, 8, 9),
    @(31, 7, 7),
    @(32, 8, 8),
    @(33, 6, 6),
    @(34, 7, 7),
    @(35, 8, 8),
    @(36, 7, 7),
    @(37, 8, 8),
    @(38, 9, 9),
    @(39, 8, 8),
    @(40, 7, 8),
    @(41, 7, 7),
    @(42, 7, 7),
    @(43, 7, 7),
    @(44, 7, 7),
    @(45, 8, 8),
    @(46, 8, 8),
    @(47, 6, 7),
    @(48, 7, 7),
    @(49, 9, 9),
    @(50, 5, 5),
    @(51, 5, 5),
    @(52, 9, 9),
    @(53, 6, 6),
    @(54, 7, 7),
    @(55, 3, 3),
    @(56, 3, 3)
)

foreach ($row in $data) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
